# Apply the content edits described by the commit:
#   - IT!B2 (the "Initial Time" year) advances from 2020 to 2021
#   - cell About!A7 loses its extra (redundant) cell style, reverting to
#     the workbook's default "Normal" style
#   - the "About" sheet becomes the active/selected sheet instead of "IT"

$wb = $excel.ActiveWorkbook

$itSheet = $wb.Worksheets.Item("IT")
$itSheet.Range("B2").Value = 2021

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A7").Style = "Normal"

# Make "About" the active sheet (was "IT").
$aboutSheet.Activate()
